# Updated symbol list on Sun Jan 15 04:08:33 UTC 2023 with GitHub Actions
#
# Refresh the crypto price/volume/hour snapshot columns (D,E,G) and swap the
# FTXToken / GateToken rows (7 & 8) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E/G hold numeric-looking text (price, %-change, hour) stored as plain text
# in the source sheet, so force NumberFormat "@" on each cell before writing the
# new value - otherwise Excel would auto-convert "298.65" / "-2.70%" into a real
# number and change the stored cell type.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '298.65'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-2.70%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '4'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.25%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '4'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.157'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.64%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '4'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07514'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.92%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '4'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.777'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.24%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '4'

# Row 7
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.677'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '7.27%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '4'

# Row 8
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.792'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.15%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '4'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9259'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.70%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '4'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1720'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.95%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '4'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07580'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.15%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '4'

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08002'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.62%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '4'

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03041'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.92%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '4'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09894'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.16%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '4'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001490'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.51%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '4'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04655'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2.34%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '4'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006249'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.86%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '4'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.456'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.57%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '4'

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.58%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '4'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3292'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.47%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '4'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1334'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.69%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '4'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.551'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '7.41%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '4'

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-4.34%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '4'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001217'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.07%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '4'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004416'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.78%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '4'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001399'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '19.72%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '4'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001807'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '8.66%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '4'

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '4'

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '4'

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '4'

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '4'

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '4'

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '4'

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '4'

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '4'

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '4'

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '4'

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '4'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01660'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.72%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '4'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04538'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.80%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '4'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006935'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-4.83%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '4'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1343'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.84%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '4'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002058'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-8.80%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '4'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01290'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-6.61%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '4'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006064'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.38%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '4'

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.95%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '4'

# Row 47
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '4'

# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '4'

# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '4'

# Row 50
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '4'

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '4'
